# Fix PER team/value assignment bug: correct the Team (col B) and
# stat (col C) values for rows 2:31 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "NOK"
$ws.Range("C2").Value = 13.03076923076923
$ws.Range("B3").Value = "POR"
$ws.Range("C3").Value = 12.85
$ws.Range("B4").Value = "NJN"
$ws.Range("C4").Value = 10.94615384615385
$ws.Range("B5").Value = "CLE"
$ws.Range("C5").Value = 13.36153846153846
$ws.Range("B6").Value = "DAL"
$ws.Range("C6").Value = 13.28125
$ws.Range("B7").Value = "ATL"
$ws.Range("C7").Value = 13.2
$ws.Range("B8").Value = "SEA"
$ws.Range("C8").Value = 12.13636363636364
$ws.Range("B9").Value = "CHA"
$ws.Range("C9").Value = 13.35
$ws.Range("B10").Value = "WAS"
$ws.Range("C10").Value = 12.82857142857143
$ws.Range("B11").Value = "MIL"
$ws.Range("C11").Value = 12.65
$ws.Range("B12").Value = "LAC"
$ws.Range("C12").Value = 11.26
$ws.Range("B13").Value = "SAS"
$ws.Range("C13").Value = 13.76
$ws.Range("B14").Value = "DET"
$ws.Range("C14").Value = 13.68461538461539
$ws.Range("B15").Value = "ORL"
$ws.Range("C15").Value = 12.71538461538461
$ws.Range("B16").Value = "UTA"
$ws.Range("C16").Value = 12.45333333333333
$ws.Range("B17").Value = "MEM"
$ws.Range("C17").Value = 13.33571428571429
$ws.Range("B18").Value = "HOU"
$ws.Range("C18").Value = 13.24166666666667
$ws.Range("B19").Value = "DEN"
$ws.Range("C19").Value = 8.784615384615385
$ws.Range("B20").Value = "LAL"
$ws.Range("C20").Value = 11.49375
$ws.Range("B21").Value = "GSW"
$ws.Range("C21").Value = 12.51333333333333
$ws.Range("B22").Value = "IND"
$ws.Range("C22").Value = 13.1
$ws.Range("B23").Value = "CHI"
$ws.Range("C23").Value = 11.125
$ws.Range("B24").Value = "PHI"
$ws.Range("C24").Value = 12.22142857142857
$ws.Range("B25").Value = "BOS"
$ws.Range("C25").Value = 12.66666666666667
$ws.Range("B26").Value = "TOR"
$ws.Range("C26").Value = 11.53076923076923
$ws.Range("B27").Value = "MIA"
$ws.Range("C27").Value = 11.8625
$ws.Range("B28").Value = "SAC"
$ws.Range("C28").Value = 14.49090909090909
$ws.Range("B29").Value = "PHO"
$ws.Range("C29").Value = 14.66153846153846
$ws.Range("B30").Value = "NYK"
$ws.Range("C30").Value = 11.77857142857143
$ws.Range("B31").Value = "MIN"
$ws.Range("C31").Value = 12.34
